$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unused rows 23:25 (records 4,5,6 removed from the bottom table)
$ws.Rows("23:25").Delete()

# Fill in Actual Results / Pass-Fail status cells that were left blank.
# Order matters for how new strings land in the shared-strings table.
$ws.Range("F21").Value = "Main section : Moble phone: 1. Photo with code is to small. 2. Main headline is sticked with down text . 3. Two phone picture is to small. 4.Headline ""Simple UI & UX"" is sticked to down text.`nTablet 1.Photo with code is to big needs to be smaller from top and left side and a bit lower 2. Headline and text section have no font color and size, and need to have space between text. 3. Phone pictures size to big. 4. Headline and text section have no font color and size, and need to have space between text.`nWeb size. 1.Photo with code is to big needs to be smaller from top and left side and a bit lower   2. Headline and text section have no font color and size, and need to have space between text. 3. Phone pictures size to big.  And need to be more to the left side.4. Headline and text section have no font color and size, and need to have space between text."

$ws.Range("F20").Value = "As expected`t`t"

$ws.Range("I20").Value = "Pass "

$ws.Range("J6").Value = "Fail"
$ws.Range("I21").Value = "Fail"

$ws.Range("F22").Value = "As expected`t`t"

$ws.Range("I22").Value = "Pass"

# Move the active selection to J7 as recorded in the edit session
$ws.Range("J7").Select()
